$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17,8).Value = 1257702.9
$ws.Cells.Item(17,10).Value = 1676437.1
$ws.Cells.Item(17,12).Value = 5029311.300000001
$ws.Cells.Item(17,14).Value = -5029647.300000001

# ALC row 58
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(58,8).Value = 2775.6667
$ws.Cells.Item(58,10).Value = 1638
$ws.Cells.Item(58,12).Value = 4914
$ws.Cells.Item(58,14).Value = -5214

# ALC row 96
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(96,8).Value = 2464.5454
$ws.Cells.Item(96,9).Value = 1118.25
$ws.Cells.Item(96,10).Value = 3233.8572
$ws.Cells.Item(96,11).Value = 3354.75
$ws.Cells.Item(96,12).Value = 9701.571599999999
$ws.Cells.Item(96,13).Value = -1981.75
$ws.Cells.Item(96,14).Value = -12447.5716

# ALC row 101
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(101,8).Value = 416.3
$ws.Cells.Item(101,9).Value = 416.3
$ws.Cells.Item(101,11).Value = 1248.9
$ws.Cells.Item(101,13).Value = 373.0999999999999

# ALC row 106
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(106,8).Value = 13335232
$ws.Cells.Item(106,9).Value = 13335232
$ws.Cells.Item(106,11).Value = 13335232
$ws.Cells.Item(106,13).Value = -13334601

# ALC row 135
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(135,8).Value = 2453.2222
$ws.Cells.Item(135,9).Value = 1635
$ws.Cells.Item(135,10).Value = 8999
$ws.Cells.Item(135,11).Value = 14715
$ws.Cells.Item(135,12).Value = 80991
$ws.Cells.Item(135,13).Value = -12180
$ws.Cells.Item(135,14).Value = -86061

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138,8).Value = 3041.0728
$ws.Cells.Item(138,9).Value = 2563.1765
$ws.Cells.Item(138,11).Value = 7689.529500000001
$ws.Cells.Item(138,13).Value = -2549.529500000001

# ALC row 141
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(141,8).Value = 3137.476
$ws.Cells.Item(141,9).Value = 3044.35
$ws.Cells.Item(141,11).Value = 9133.049999999999
$ws.Cells.Item(141,13).Value = -3953.049999999999

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32,8).Value = 8505652
$ws.Cells.Item(32,9).Value = 4762988.5
$ws.Cells.Item(32,10).Value = 17862312
$ws.Cells.Item(32,11).Value = 4762988.5
$ws.Cells.Item(32,12).Value = 17862312
$ws.Cells.Item(32,13).Value = -4762701.5
$ws.Cells.Item(32,14).Value = -17862886

# ARM row 45
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45,8).Value = 2783
$ws.Cells.Item(45,9).Value = 2585.3845
$ws.Cells.Item(45,11).Value = 2585.3845
$ws.Cells.Item(45,13).Value = -2208.3845

# ARM row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61,8).Value = 2109.6135
$ws.Cells.Item(61,9).Value = 1876.85
$ws.Cells.Item(61,10).Value = 4437.25
$ws.Cells.Item(61,11).Value = 1876.85
$ws.Cells.Item(61,12).Value = 4437.25
$ws.Cells.Item(61,13).Value = -1664.85
$ws.Cells.Item(61,14).Value = -4861.25

# ARM row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74,8).Value = 1687.3077
$ws.Cells.Item(74,9).Value = 1269.6666
$ws.Cells.Item(74,11).Value = 1269.6666
$ws.Cells.Item(74,13).Value = -395.6666

# ARM row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77,8).Value = 1687.3077
$ws.Cells.Item(77,9).Value = 1269.6666
$ws.Cells.Item(77,11).Value = 6348.333000000001
$ws.Cells.Item(77,13).Value = -1980.333000000001

# ARM row 122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122,8).Value = 3672.0132
$ws.Cells.Item(122,9).Value = 2662.288
$ws.Cells.Item(122,10).Value = 7176.353
$ws.Cells.Item(122,11).Value = 7986.864
$ws.Cells.Item(122,12).Value = 21529.059
$ws.Cells.Item(122,13).Value = -5536.864
$ws.Cells.Item(122,14).Value = -26429.059

# ARM row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132,8).Value = 3406.3215
$ws.Cells.Item(132,9).Value = 3116.158
$ws.Cells.Item(132,11).Value = 9348.474
$ws.Cells.Item(132,13).Value = -6818.474

# ARM row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136,8).Value = 2109.6135
$ws.Cells.Item(136,9).Value = 1876.85
$ws.Cells.Item(136,10).Value = 4437.25
$ws.Cells.Item(136,11).Value = 5630.549999999999
$ws.Cells.Item(136,12).Value = 13311.75
$ws.Cells.Item(136,13).Value = -3080.549999999999
$ws.Cells.Item(136,14).Value = -18411.75

# ARM row 139
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(139,8).Value = 0
$ws.Cells.Item(139,10).Value = 0
$ws.Cells.Item(139,12).Value = 0
$ws.Cells.Item(139,14).ClearContents()

# BSM row 82
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(82,8).Value = 4299
$ws.Cells.Item(82,9).Value = 4299
$ws.Cells.Item(82,10).Value = 0
$ws.Cells.Item(82,11).Value = 4299
$ws.Cells.Item(82,12).Value = 0
$ws.Cells.Item(82,13).Value = -3916
$ws.Cells.Item(82,14).ClearContents()

# BSM row 85
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(85,8).Value = 4299
$ws.Cells.Item(85,9).Value = 4299
$ws.Cells.Item(85,10).Value = 0
$ws.Cells.Item(85,11).Value = 4299
$ws.Cells.Item(85,12).Value = 0
$ws.Cells.Item(85,13).Value = -2973
$ws.Cells.Item(85,14).ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(105,8).Value = 4174.9
$ws.Cells.Item(105,9).Value = 0
$ws.Cells.Item(105,10).Value = 4174.9
$ws.Cells.Item(105,11).Value = 0
$ws.Cells.Item(105,12).Value = 4174.9
$ws.Cells.Item(105,13).ClearContents()
$ws.Cells.Item(105,14).Value = -7668.9

# BSM row 107
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(107,8).Value = 1674
$ws.Cells.Item(107,9).Value = 1425
$ws.Cells.Item(107,11).Value = 1425
$ws.Cells.Item(107,13).Value = 495

# CRP row 6
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6,8).Value = 6701980
$ws.Cells.Item(6,9).Value = 6701980
$ws.Cells.Item(6,11).Value = 6701980
$ws.Cells.Item(6,13).Value = -6701867

# CRP row 94
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(94,8).Value = 3180.1428
$ws.Cells.Item(94,9).Value = 2949.5
$ws.Cells.Item(94,10).Value = 3218.5833
$ws.Cells.Item(94,11).Value = 2949.5
$ws.Cells.Item(94,12).Value = 3218.5833
$ws.Cells.Item(94,13).Value = -2498.5
$ws.Cells.Item(94,14).Value = -4120.5833

# CRP row 111
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(111,8).Value = 48996.332
$ws.Cells.Item(111,10).Value = 48996.332
$ws.Cells.Item(111,12).Value = 48996.332
$ws.Cells.Item(111,14).Value = -57176.332

# CRP row 118
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(118,8).Value = 116999
$ws.Cells.Item(118,10).Value = 116999
$ws.Cells.Item(118,12).Value = 116999
$ws.Cells.Item(118,14).Value = -120313

# CRP row 132
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132,8).Value = 2746.4614
$ws.Cells.Item(132,9).Value = 2870.4
$ws.Cells.Item(132,11).Value = 8611.200000000001
$ws.Cells.Item(132,13).Value = -6081.200000000001

# CRP row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134,8).Value = 2112.361
$ws.Cells.Item(134,9).Value = 1564.5667
$ws.Cells.Item(134,11).Value = 4693.7001
$ws.Cells.Item(134,13).Value = -2158.7001

# CRP row 140
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(140,8).Value = 999999
$ws.Cells.Item(140,10).Value = 999999
$ws.Cells.Item(140,12).Value = 999999
$ws.Cells.Item(140,14).Value = -1010359

# CUL row 4
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4,8).Value = 86045490
$ws.Cells.Item(4,9).Value = 66844696
$ws.Cells.Item(4,11).Value = 200534088
$ws.Cells.Item(4,13).Value = -200533976

# CUL row 38
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(38,8).Value = 618.75
$ws.Cells.Item(38,10).Value = 851.7778
$ws.Cells.Item(38,12).Value = 2555.3334
$ws.Cells.Item(38,14).Value = -3249.3334

# CUL row 46
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(46,8).Value = 231.72728
$ws.Cells.Item(46,9).Value = 232.4
$ws.Cells.Item(46,11).Value = 697.2
$ws.Cells.Item(46,13).Value = -606.2

# CUL row 107
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(107,8).Value = 631.94446
$ws.Cells.Item(107,10).Value = 525.06665
$ws.Cells.Item(107,12).Value = 1575.19995
$ws.Cells.Item(107,14).Value = -5415.19995

# CUL row 122
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(122,8).Value = 861
$ws.Cells.Item(122,9).Value = 895.6667
$ws.Cells.Item(122,10).Value = 809
$ws.Cells.Item(122,11).Value = 8061.0003
$ws.Cells.Item(122,12).Value = 7281
$ws.Cells.Item(122,13).Value = -5611.0003
$ws.Cells.Item(122,14).Value = -12181

# CUL row 131
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(131,8).Value = 1704.5358
$ws.Cells.Item(131,9).Value = 1282.3334
$ws.Cells.Item(131,10).Value = 1819.6818
$ws.Cells.Item(131,11).Value = 3847.0002
$ws.Cells.Item(131,12).Value = 5459.0454
$ws.Cells.Item(131,13).Value = 1192.9998
$ws.Cells.Item(131,14).Value = -15539.0454

# CUL row 140
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(140,8).Value = 22224694
$ws.Cells.Item(140,9).Value = 22224694
$ws.Cells.Item(140,11).Value = 66674082
$ws.Cells.Item(140,13).Value = -66668902

# GSM row 33
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(33,8).Value = 23266
$ws.Cells.Item(33,10).Value = 24900
$ws.Cells.Item(33,12).Value = 24900
$ws.Cells.Item(33,14).Value = -25404

# GSM row 123
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(123,8).Value = 34516.75
$ws.Cells.Item(123,10).Value = 34516.75
$ws.Cells.Item(123,12).Value = 34516.75
$ws.Cells.Item(123,14).Value = -39416.75

# GSM row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132,8).Value = 4270.7
$ws.Cells.Item(132,9).Value = 4270.7
$ws.Cells.Item(132,11).Value = 12812.1
$ws.Cells.Item(132,13).Value = -10282.1

# LTW row 61
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61,8).Value = 1537
$ws.Cells.Item(61,9).Value = 1470.8572
$ws.Cells.Item(61,11).Value = 1470.8572
$ws.Cells.Item(61,13).Value = -1268.8572

# LTW row 82
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(82,8).Value = 2297.84
$ws.Cells.Item(82,10).Value = 2932.9285
$ws.Cells.Item(82,12).Value = 2932.9285
$ws.Cells.Item(82,14).Value = -3654.9285

# LTW row 85
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(85,8).Value = 2297.84
$ws.Cells.Item(85,10).Value = 2932.9285
$ws.Cells.Item(85,12).Value = 2932.9285
$ws.Cells.Item(85,14).Value = -5428.9285

# LTW row 100
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(100,8).Value = 2313
$ws.Cells.Item(100,9).Value = 2313
$ws.Cells.Item(100,11).Value = 2313
$ws.Cells.Item(100,13).Value = -1772

# LTW row 113
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113,8).Value = 1537
$ws.Cells.Item(113,9).Value = 1470.8572
$ws.Cells.Item(113,11).Value = 1470.8572
$ws.Cells.Item(113,13).Value = 699.1428000000001

# LTW row 132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132,8).Value = 3621.889
$ws.Cells.Item(132,9).Value = 2738.4
$ws.Cells.Item(132,10).Value = 4726.25
$ws.Cells.Item(132,11).Value = 8215.200000000001
$ws.Cells.Item(132,12).Value = 14178.75
$ws.Cells.Item(132,13).Value = -5685.200000000001
$ws.Cells.Item(132,14).Value = -19238.75

# LTW row 135
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(135,8).Value = 49979.5
$ws.Cells.Item(135,10).Value = 49979.5
$ws.Cells.Item(135,12).Value = 49979.5
$ws.Cells.Item(135,14).Value = -60119.5

# LTW row 136
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136,8).Value = 2962.5217
$ws.Cells.Item(136,9).Value = 2306.1177
$ws.Cells.Item(136,10).Value = 4822.3335
$ws.Cells.Item(136,11).Value = 6918.353099999999
$ws.Cells.Item(136,12).Value = 14467.0005
$ws.Cells.Item(136,13).Value = -4368.353099999999
$ws.Cells.Item(136,14).Value = -19567.0005

# WVR row 39
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(39,8).Value = 29458.5
$ws.Cells.Item(39,9).Value = 29458.5
$ws.Cells.Item(39,11).Value = 29458.5
$ws.Cells.Item(39,13).Value = -29045.5

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132,8).Value = 1576.7222
$ws.Cells.Item(132,10).Value = 2309.1667
$ws.Cells.Item(132,12).Value = 6927.500100000001
$ws.Cells.Item(132,14).Value = -11987.5001
